$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newValue) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $newValue
    $c.ClearFormats()
}

Set-TextValue "D2" "51.773.90"
Set-TextValue "E2" "  +1.64%  "
Set-TextValue "D3" "3.035.71"
Set-TextValue "E3" "  +3.05%  "
Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.04%  "
Set-TextValue "D5" "380.76"
Set-TextValue "E5" "  +0.46%  "
Set-TextValue "D6" "103.54"
Set-TextValue "E6" "  +2.42%  "
Set-TextValue "E7" "  +0.93%  "
Set-TextValue "E8" "  -0.01%  "
Set-TextValue "D9" "0.597"
Set-TextValue "E9" "  +2.82%  "
Set-TextValue "D10" "37.15"
Set-TextValue "E10" "  +2.91%  "
Set-TextValue "E11" "  -0.19%  "
Set-TextValue "D12" "0.0863"
Set-TextValue "E12" "  +1.86%  "
Set-TextValue "D13" "3.524.79"
Set-TextValue "E13" "  +3.56%  "
Set-TextValue "D14" "18.62"
Set-TextValue "E14" "  +1.90%  "
Set-TextValue "D15" "7.77"
Set-TextValue "E15" "  +0.11%  "
Set-TextValue "D16" "3.057.10"
Set-TextValue "E16" "  +3.79%  "
Set-TextValue "D17" "0.982"
Set-TextValue "E17" "  -1.79%  "
Set-TextValue "D18" "10.55"
Set-TextValue "E18" "  -12.65%  "
Set-TextValue "D19" "51.772.11"
Set-TextValue "E19" "  +1.66%  "
Set-TextValue "D20" "3.06"
Set-TextValue "E20" "  -0.40%  "
Set-TextValue "D21" "12.57"
Set-TextValue "E21" "  +1.55%  "
Set-TextValue "D22" "0.0₃0966"
Set-TextValue "E22" "  +1.69%  "
Set-TextValue "D23" "70.09"
Set-TextValue "E23" "  +0.99%  "
Set-TextValue "D24" "269.28"
Set-TextValue "E24" "  +1.16%  "
Set-TextValue "D25" "3.18"
Set-TextValue "E25" "  -0.58%  "
Set-TextValue "E26" "  +0.67%  "
Set-TextValue "D27" "7.56"
Set-TextValue "E27" "  +7.19%  "
Set-TextValue "E28" "  +6.34%  "
Set-TextValue "D29" "26.34"
Set-TextValue "E29" "  +3.01%  "
Set-TextValue "E30" "  -0.05%  "
Set-TextValue "E31" "  +0.81%  "
Set-TextValue "D32" "10.33"
Set-TextValue "E32" "  +2.26%  "
Set-TextValue "D33" "34.32"
Set-TextValue "E33" "  +2.72%  "
Set-TextValue "E34" "  +0.09%  "
Set-TextValue "D35" "50.49"
Set-TextValue "E35" "  +0.12%  "
Set-TextValue "D36" "0.0452"
Set-TextValue "E36" "  +5.02%  "
Set-TextValue "E37" "  -0.09%  "
Set-TextValue "E38" "  +8.82%  "
Set-TextValue "E39" "  +11.36%  "
Set-TextValue "D40" "17.18"
Set-TextValue "E40" "  +3.69%  "
Set-TextValue "D41" "1.87"
Set-TextValue "E41" "  +3.70%  "
Set-TextValue "E42" "  +3.37%  "
Set-TextValue "E43" "  +0.39%  "
Set-TextValue "D44" "127.49"
Set-TextValue "E44" "  +7.14%  "
Set-TextValue "D45" "3.77"
Set-TextValue "E45" "  +7.14%  "
Set-TextValue "D46" "21.94"
Set-TextValue "E46" "  +2.90%  "
Set-TextValue "E47" "  +6.02%  "
Set-TextValue "D48" "2.41"
Set-TextValue "E48" "  +3.87%  "
Set-TextValue "D49" "2.037.16"
Set-TextValue "E49" "  +1.81%  "
Set-TextValue "D50" "3.338.59"
Set-TextValue "E50" "  +3.04%  "
Set-TextValue "E51" "  +2.56%  "
